# Revert "Merge branch 'image_file_formats'"
# - restore the cached datetimeFigureOut field text on the slide master,
#   the slide layouts and the notes master back to "15.01.21"
# - drop the 3rd slide that the merge had added

$p = $ppt.ActivePresentation

# --- slide master ------------------------------------------------------
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $sh = $masterShapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.TextFrame.TextRange.Text -eq "04.01.22") {
            $sh.TextFrame.TextRange.Text = "15.01.21"
        }
    }
}

# --- every slide layout --------------------------------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $shapes = $layouts.Item($i).Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq "04.01.22") {
                $sh.TextFrame.TextRange.Text = "15.01.21"
            }
        }
    }
}

# --- notes master ----------------------------------------------------------
$notesShapes = $p.NotesMaster.Shapes
for ($i = 1; $i -le $notesShapes.Count; $i++) {
    $sh = $notesShapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.TextFrame.TextRange.Text -eq "04.01.22") {
            $sh.TextFrame.TextRange.Text = "15.01.21"
        }
    }
}

# --- remove the slide that the merge had introduced -------------------------
$p.Slides.Item(3).Delete()
